# changed factor as rolling mean / rolling std
#
# - remove the "scale" calibration parameter (label + value), which shifts
#   the "window" row up from row 7 to row 6
# - update end_date (B3) to the new date
# - update startPrice (B4) to the new rolling-mean-based value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "scale" row (row 6) - shifts "window" (old row 7) up to row 6
$ws.Rows.Item(6).Delete()

# Update end_date value (row 3) - force text so it isn't auto-parsed as a date
# serial number, then clear the temporary number format so the cell keeps the
# default (unstyled) appearance, matching the rest of the column.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2021-12-14"
$ws.Range("B3").ClearFormats()

# Update startPrice value (row 4)
$ws.Range("B4").Value = 4668.97021484375
